$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8
$ws.Range("A8").Value = "electricity, high voltage"
$ws.Range("B8").Value = "treatment of bagasse, from sweet sorghum, in heat and power co-generation unit, 6400kW thermal"
$ws.Range("C8").Value = "GLO"
$ws.Range("D8").Value = "ecoinvent38_cutoff"

# Row 9
$ws.Range("A9").Value = "heat, district or industrial, other than natural gas"
$ws.Range("B9").Value = "treatment of bagasse, from sweet sorghum, in heat and power co-generation unit, 6400kW thermal"
$ws.Range("C9").Value = "GLO"
$ws.Range("D9").Value = "ecoinvent38_cutoff"

# Row 10
$ws.Range("A10").Value = "electricity, low voltage"
$ws.Range("B10").Value = "wood pellets, burned in stirling heat and power co-generation unit, 3kW electrical, future"
$ws.Range("C10").Value = "CH"
$ws.Range("D10").Value = "ecoinvent38_cutoff"

# Row 11
$ws.Range("A11").Value = "heat, future"
$ws.Range("B11").Value = "wood pellets, burned in stirling heat and power co-generation unit, 3kW electrical, future"
$ws.Range("C11").Value = "CH"
$ws.Range("D11").Value = "ecoinvent38_cutoff"

# Update selection to match target (A12 selected, as if user moved to next empty row)
$ws.Range("A12").Select()
